$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1107.4667
$ws.Range("I19").Value = 1041.1428
$ws.Range("J19").Value = 1165.5
$ws.Range("K19").Value = 1041.1428
$ws.Range("L19").Value = 1165.5
$ws.Range("M19").Value = -866.1428000000001
$ws.Range("N19").Value = -1515.5
$ws.Range("H70").Value = 1133
$ws.Range("J70").Value = 839.44446
$ws.Range("L70").Value = 2518.33338
$ws.Range("N70").Value = -3058.33338
$ws.Range("H73").Value = 1133
$ws.Range("J73").Value = 839.44446
$ws.Range("L73").Value = 2518.33338
$ws.Range("N73").Value = -4390.33338
$ws.Range("H103").Value = 564.2857
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 491.66666
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 1474.99998
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -2646.99998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 79170.53999999999
$ws.Range("I102").Value = 112987.664
$ws.Range("J102").Value = 3082
$ws.Range("K102").Value = 112987.664
$ws.Range("L102").Value = 3082
$ws.Range("M102").Value = -111365.664
$ws.Range("N102").Value = -6326
$ws.Range("H122").Value = 1844.5
$ws.Range("I122").Value = 2790
$ws.Range("J122").Value = 899
$ws.Range("K122").Value = 8370
$ws.Range("L122").Value = 2697
$ws.Range("M122").Value = -5920
$ws.Range("N122").Value = -7597

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 517.85
$ws.Range("I94").Value = 483.2
$ws.Range("J94").Value = 621.8
$ws.Range("K94").Value = 483.2
$ws.Range("L94").Value = 621.8
$ws.Range("M94").Value = -32.19999999999999
$ws.Range("N94").Value = -1523.8
$ws.Range("H99").Value = 1924.375
$ws.Range("I99").Value = 1482.8572
$ws.Range("J99").Value = 2267.7778
$ws.Range("K99").Value = 1482.8572
$ws.Range("L99").Value = 2267.7778
$ws.Range("M99").Value = 15.14280000000008
$ws.Range("N99").Value = -5263.7778
$ws.Range("H107").Value = 62529280
$ws.Range("I107").Value = 250114290
$ws.Range("J107").Value = 944.4167
$ws.Range("K107").Value = 250114290
$ws.Range("L107").Value = 944.4167
$ws.Range("M107").Value = -250112370
$ws.Range("N107").Value = -4784.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1430.579
$ws.Range("I58").Value = 1192.1765
$ws.Range("J58").Value = 3457
$ws.Range("K58").Value = 1192.1765
$ws.Range("L58").Value = 3457
$ws.Range("M58").Value = -989.1765
$ws.Range("N58").Value = -3863
$ws.Range("H99").Value = 14410.625
$ws.Range("I99").Value = 2624
$ws.Range("J99").Value = 21482.6
$ws.Range("K99").Value = 2624
$ws.Range("L99").Value = 21482.6
$ws.Range("M99").Value = -1126
$ws.Range("N99").Value = -24478.6
$ws.Range("H126").Value = 14410.625
$ws.Range("I126").Value = 2624
$ws.Range("J126").Value = 21482.6
$ws.Range("K126").Value = 7872
$ws.Range("L126").Value = 64447.8
$ws.Range("M126").Value = -5402
$ws.Range("N126").Value = -69387.79999999999
$ws.Range("H132").Value = 88241660
$ws.Range("I132").Value = 90915896
$ws.Range("J132").Value = 83338920
$ws.Range("K132").Value = 272747688
$ws.Range("L132").Value = 250016760
$ws.Range("M132").Value = -272745158
$ws.Range("N132").Value = -250021820
$ws.Range("H136").Value = 1430.579
$ws.Range("I136").Value = 1192.1765
$ws.Range("J136").Value = 3457
$ws.Range("K136").Value = 3576.5295
$ws.Range("L136").Value = 10371
$ws.Range("M136").Value = -1026.5295
$ws.Range("N136").Value = -15471

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 665.9286
$ws.Range("I113").Value = 700
$ws.Range("J113").Value = 652.3
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 1956.9
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -6296.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3793.611
$ws.Range("I102").Value = 3286.5
$ws.Range("J102").Value = 4199.3
$ws.Range("K102").Value = 3286.5
$ws.Range("L102").Value = 4199.3
$ws.Range("M102").Value = -1664.5
$ws.Range("N102").Value = -7443.3
$ws.Range("H113").Value = 1433.25
$ws.Range("I113").Value = 729
$ws.Range("J113").Value = 1785.375
$ws.Range("K113").Value = 729
$ws.Range("L113").Value = 1785.375
$ws.Range("M113").Value = 1441
$ws.Range("N113").Value = -6125.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 51833.9
$ws.Range("I40").Value = 126150.5
$ws.Range("J40").Value = 2289.5
$ws.Range("K40").Value = 126150.5
$ws.Range("L40").Value = 2289.5
$ws.Range("M40").Value = -126014.5
$ws.Range("N40").Value = -2561.5
$ws.Range("H61").Value = 2046.2142
$ws.Range("I61").Value = 2240.6
$ws.Range("J61").Value = 1938.2222
$ws.Range("K61").Value = 2240.6
$ws.Range("L61").Value = 1938.2222
$ws.Range("M61").Value = -2038.6
$ws.Range("N61").Value = -2342.2222
$ws.Range("H68").Value = 3373.4285
$ws.Range("I68").Value = 2446.6667
$ws.Range("J68").Value = 3626.182
$ws.Range("K68").Value = 2446.6667
$ws.Range("L68").Value = 3626.182
$ws.Range("M68").Value = -1697.6667
$ws.Range("N68").Value = -5124.182
$ws.Range("H71").Value = 3373.4285
$ws.Range("I71").Value = 2446.6667
$ws.Range("J71").Value = 3626.182
$ws.Range("K71").Value = 12233.3335
$ws.Range("L71").Value = 18130.91
$ws.Range("M71").Value = -8489.333500000001
$ws.Range("N71").Value = -25618.91
$ws.Range("H82").Value = 1865.8334
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 2039
$ws.Range("K82").Value = 1000
$ws.Range("L82").Value = 2039
$ws.Range("M82").Value = -639
$ws.Range("N82").Value = -2761
$ws.Range("H85").Value = 1865.8334
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 2039
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 2039
$ws.Range("M85").Value = 248
$ws.Range("N85").Value = -4535
$ws.Range("H93").Value = 2126.4211
$ws.Range("J93").Value = 1356.5555
$ws.Range("L93").Value = 1356.5555
$ws.Range("N93").Value = -3852.5555
$ws.Range("H100").Value = 2472.5
$ws.Range("I100").Value = 2075
$ws.Range("J100").Value = 2870
$ws.Range("K100").Value = 2075
$ws.Range("L100").Value = 2870
$ws.Range("M100").Value = -1534
$ws.Range("N100").Value = -3952
$ws.Range("H113").Value = 2046.2142
$ws.Range("I113").Value = 2240.6
$ws.Range("J113").Value = 1938.2222
$ws.Range("K113").Value = 2240.6
$ws.Range("L113").Value = 1938.2222
$ws.Range("M113").Value = -70.59999999999991
$ws.Range("N113").Value = -6278.2222
$ws.Range("H122").Value = 3022.111
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3066.3333
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9198.999899999999
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -14098.9999
$ws.Range("H132").Value = 2974.85
$ws.Range("I132").Value = 2942.606
$ws.Range("J132").Value = 3126.8572
$ws.Range("K132").Value = 8827.818000000001
$ws.Range("L132").Value = 9380.571599999999
$ws.Range("M132").Value = -6297.818000000001
$ws.Range("N132").Value = -14440.5716
$ws.Range("H133").Value = 33265.89
$ws.Range("J133").Value = 33265.89
$ws.Range("L133").Value = 33265.89
$ws.Range("N133").Value = -38325.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 200001950
$ws.Range("I96").Value = 250002100
$ws.Range("J96").Value = 1400
$ws.Range("K96").Value = 250002100
$ws.Range("L96").Value = 1400
$ws.Range("M96").Value = -250000727
$ws.Range("N96").Value = -4146
$ws.Range("H122").Value = 1998.3334
$ws.Range("I122").Value = 1622.75
$ws.Range("J122").Value = 2749.5
$ws.Range("K122").Value = 4868.25
$ws.Range("L122").Value = 8248.5
$ws.Range("M122").Value = -2418.25
$ws.Range("N122").Value = -13148.5
$ws.Range("H136").Value = 2284
$ws.Range("I136").Value = 614.5833
$ws.Range("J136").Value = 4788.125
$ws.Range("K136").Value = 1843.7499
$ws.Range("L136").Value = 14364.375
$ws.Range("M136").Value = 706.2501
$ws.Range("N136").Value = -19464.375
